# Refresh the cryptocurrency price/volume table (cryptos.xlsx, Sheet1).
# A new coin (BitDAO) was inserted at row 23, shifting most of the
# subsequent rows down by one (rows pinned to Aave/RenderToken/Aptos
# keep their own coin). Prices (col D) and 1h volume deltas (col E) were
# also refreshed for every row. All of these columns are plain text in
# the sheet (e.g. "30.509.23", "  -1.45%  "), so numeric-looking values
# are written with a temporary text number-format to stop Excel from
# re-interpreting them (and silently dropping trailing zeros), then the
# format is cleared again so no stray cell style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.528.82"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "1.913.31"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4762"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2846"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06680"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.80"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "101.28"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.55%  "
$ws.Range("D12").Value = "1.915.35"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07683"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.230"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6706"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.19%  "
$ws.Range("D16").Value = "30.527.41"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "255.87"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -8.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007489"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.67"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.402"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("B23").Value = "BitDAO"
$ws.Range("C23").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.4507"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -9.95%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.304"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.349"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.00%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.45%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.058"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.17%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.713"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.93%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1008"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.374"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.513"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.06%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.263"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04724"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7287"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.111"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.56%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9994"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.710"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01913"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.33%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.610"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.98"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.02%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.224"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.962"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.89%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8614"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.66%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.17"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.64%  "
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4244"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.44%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9997"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.416"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.87%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "984.72"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1198"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.04%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.81"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.20%  "
